# Case_2_216/res_line/loading_percent: update simulated loading-percent
# results for the 380 kV case (commit: "case with 380 kV done").
# Rewrites the per-timestep loading-percent values (rows 2-25, i.e. time
# steps 0-23) for the non-zero line columns (B, C, E, F, G, I, K, L, N);
# columns that are always 0 (D, H, J, M, O) and the time-step index in
# column A are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 17.20634126269098
$ws.Range("C2").Value = 4.797472585559071
$ws.Range("E2").Value = 9.500935766553042
$ws.Range("F2").Value = 16.86991607391233
$ws.Range("G2").Value = 3.734061416682816
$ws.Range("I2").Value = 32.68088263532354
$ws.Range("K2").Value = 15.14775986030986
$ws.Range("L2").Value = 10.51751424675212
$ws.Range("N2").Value = 22.94409605711508
$ws.Range("B3").Value = 17.04669131259605
$ws.Range("C3").Value = 4.649641312072581
$ws.Range("E3").Value = 9.504947920102081
$ws.Range("F3").Value = 15.89584955866808
$ws.Range("G3").Value = 3.736888249148664
$ws.Range("I3").Value = 32.73230220484155
$ws.Range("K3").Value = 15.04564818920103
$ws.Range("L3").Value = 10.51025438244832
$ws.Range("N3").Value = 23.00387421685611
$ws.Range("B4").Value = 16.95236475760775
$ws.Range("C4").Value = 4.55808695872246
$ws.Range("E4").Value = 9.508874970841052
$ws.Range("F4").Value = 15.26997757108491
$ws.Range("G4").Value = 3.738714469603861
$ws.Range("I4").Value = 32.76952747332891
$ws.Range("K4").Value = 14.98633375597946
$ws.Range("L4").Value = 10.50773741404935
$ws.Range("N4").Value = 23.04254969889868
$ws.Range("B5").Value = 16.9148964878573
$ws.Range("C5").Value = 4.520646380807237
$ws.Range("E5").Value = 9.510843843397264
$ws.Range("F5").Value = 15.00819731993403
$ws.Range("G5").Value = 3.739481513601824
$ws.Range("I5").Value = 32.78611582951385
$ws.Range("K5").Value = 14.96303486658417
$ws.Range("L5").Value = 10.50720093684365
$ws.Range("N5").Value = 23.05880668294295
$ws.Range("B6").Value = 16.90873467160718
$ws.Range("C6").Value = 4.514423468050415
$ws.Range("E6").Value = 9.511193047496022
$ws.Range("F6").Value = 14.96433081551593
$ws.Range("G6").Value = 3.739610262732525
$ws.Range("I6").Value = 32.78895591489899
$ws.Range("K6").Value = 14.95921938692506
$ws.Range("L6").Value = 10.50714143263944
$ws.Range("N6").Value = 23.06153614483782
$ws.Range("B7").Value = 16.95185546576045
$ws.Range("C7").Value = 4.557582464070607
$ws.Range("E7").Value = 9.508900030785107
$ws.Range("F7").Value = 15.26647399323137
$ws.Range("G7").Value = 3.738724721627222
$ws.Range("I7").Value = 32.76974544922793
$ws.Range("K7").Value = 14.98601597995362
$ws.Range("L7").Value = 10.50772819673131
$ws.Range("N7").Value = 23.04276693545912
$ws.Range("B8").Value = 17.15055100874783
$ws.Range("C8").Value = 4.746701764345418
$ws.Range("E8").Value = 9.502015745678825
$ws.Range("F8").Value = 16.53996406344768
$ws.Range("G8").Value = 3.735017365762915
$ws.Range("I8").Value = 32.69743708362248
$ws.Range("K8").Value = 15.11186315365264
$ws.Range("L8").Value = 10.51460908147497
$ws.Range("N8").Value = 22.96429879803712
$ws.Range("B9").Value = 17.56768145691151
$ws.Range("C9").Value = 5.108639460000029
$ws.Range("E9").Value = 9.500101506645008
$ws.Range("F9").Value = 19.00274580682531
$ws.Range("G9").Value = 3.728461998372236
$ws.Range("I9").Value = 32.60061843079068
$ws.Range("K9").Value = 15.38446768175655
$ws.Range("L9").Value = 10.54343642941255
$ws.Range("N9").Value = 22.82603421739135
$ws.Range("B10").Value = 17.88830870162573
$ws.Range("C10").Value = 5.365822948130824
$ws.Range("E10").Value = 9.505717128372373
$ws.Range("F10").Value = 20.67494806633232
$ws.Range("G10").Value = 3.72407645712029
$ws.Range("I10").Value = 32.55705936853931
$ws.Range("K10").Value = 15.59903418584196
$ws.Range("L10").Value = 10.57386313193071
$ws.Range("N10").Value = 22.73392348189002
$ws.Range("B11").Value = 18.03667182410747
$ws.Range("C11").Value = 5.480323473282117
$ws.Range("E11").Value = 9.509785610412704
$ws.Range("F11").Value = 21.3917225636224
$ws.Range("G11").Value = 3.722173805897199
$ws.Range("I11").Value = 32.5432582708651
$ws.Range("K11").Value = 15.69941505919248
$ws.Range("L11").Value = 10.58968585543347
$ws.Range("N11").Value = 22.69406749970372
$ws.Range("B12").Value = 18.09316487341578
$ws.Range("C12").Value = 5.523279356270188
$ws.Range("E12").Value = 9.511542786412567
$ws.Range("F12").Value = 21.65686569030329
$ws.Range("G12").Value = 3.721466519799878
$ws.Range("I12").Value = 32.53889869434875
$ws.Range("K12").Value = 15.73779628424933
$ws.Range("L12").Value = 10.59595965334796
$ws.Range("N12").Value = 22.67926863378396
$ws.Range("B13").Value = 18.08098503503335
$ws.Range("C13").Value = 5.514046620322814
$ws.Range("E13").Value = 9.51115473675323
$ws.Range("F13").Value = 21.60004134736742
$ws.Range("G13").Value = 3.721618260315867
$ws.Range("I13").Value = 32.5397990399311
$ws.Range("K13").Value = 15.7295142605497
$ws.Range("L13").Value = 10.59459598173239
$ws.Range("N13").Value = 22.68244277727421
$ws.Range("B14").Value = 18.04131358921584
$ws.Range("C14").Value = 5.483865794348175
$ws.Range("E14").Value = 9.509925841866131
$ws.Range("F14").Value = 21.4136618050453
$ws.Range("G14").Value = 3.722115352784123
$ws.Range("I14").Value = 32.54288222887591
$ws.Range("K14").Value = 15.70256547906694
$ws.Range("L14").Value = 10.59019636868529
$ws.Range("N14").Value = 22.69284410371778
$ws.Range("B15").Value = 18.01705273333585
$ws.Range("C15").Value = 5.465325445863559
$ws.Range("E15").Value = 9.509201270184631
$ws.Range("F15").Value = 21.29868154950795
$ws.Range("G15").Value = 3.722421553987551
$ws.Range("I15").Value = 32.54488367534934
$ws.Range("K15").Value = 15.68610575796137
$ws.Range("L15").Value = 10.58753812377583
$ws.Range("N15").Value = 22.69925345630537
$ws.Range("B16").Value = 17.87865906396031
$ws.Range("C16").Value = 5.358285993106724
$ws.Range("E16").Value = 9.505481610083146
$ws.Range("F16").Value = 20.62722412089977
$ws.Range("G16").Value = 3.724202651930217
$ws.Range("I16").Value = 32.55808248122742
$ws.Range("K16").Value = 15.59252734887475
$ws.Range("L16").Value = 10.57286871404911
$ws.Range("N16").Value = 22.73656929028224
$ws.Range("B17").Value = 17.79436626043908
$ws.Range("C17").Value = 5.291949056237369
$ws.Range("E17").Value = 9.503586674058392
$ws.Range("F17").Value = 20.20408069597325
$ws.Range("G17").Value = 3.725318899195649
$ws.Range("I17").Value = 32.56772124588272
$ws.Range("K17").Value = 15.53580949285453
$ws.Range("L17").Value = 10.56437506302864
$ws.Range("N17").Value = 22.75998497119373
$ws.Range("B18").Value = 17.74612233415838
$ws.Range("C18").Value = 5.253561462840419
$ws.Range("E18").Value = 9.50263933059685
$ws.Range("F18").Value = 19.95656407809801
$ws.Range("G18").Value = 3.725969632219943
$ws.Range("I18").Value = 32.57383120047819
$ws.Range("K18").Value = 15.50344975814172
$ws.Range("L18").Value = 10.55967643308741
$ws.Range("N18").Value = 22.77364562235558
$ws.Range("B19").Value = 17.72983033545265
$ws.Range("C19").Value = 5.240525606116127
$ws.Range("E19").Value = 9.502343097158063
$ws.Range("F19").Value = 19.87204792380568
$ws.Range("G19").Value = 3.726191455233952
$ws.Range("I19").Value = 32.57599707681263
$ws.Range("K19").Value = 15.49253935018894
$ws.Range("L19").Value = 10.55811770497786
$ws.Range("N19").Value = 22.77830397196986
$ws.Range("B20").Value = 17.80331498922034
$ws.Range("C20").Value = 5.299035119021853
$ws.Range("E20").Value = 9.50377364548779
$ws.Range("F20").Value = 20.24955283636154
$ws.Range("G20").Value = 3.725199173183591
$ws.Range("I20").Value = 32.5666365909733
$ws.Range("K20").Value = 15.54182021157974
$ws.Range("L20").Value = 10.56525992441383
$ws.Range("N20").Value = 22.75747240479805
$ws.Range("B21").Value = 18.05295800491938
$ws.Range("C21").Value = 5.492741902936666
$ws.Range("E21").Value = 9.510280931417077
$ws.Range("F21").Value = 21.46857628470577
$ws.Range("G21").Value = 3.721968986837671
$ws.Range("I21").Value = 32.54195308967162
$ws.Range("K21").Value = 15.71047120973139
$ws.Range("L21").Value = 10.59148101035983
$ws.Range("N21").Value = 22.68978101360376
$ws.Range("B22").Value = 18.2179048280769
$ws.Range("C22").Value = 5.616974138922862
$ws.Range("E22").Value = 9.51579529642709
$ws.Range("F22").Value = 22.22866616901552
$ws.Range("G22").Value = 3.719934817790531
$ws.Range("I22").Value = 32.53087280379226
$ws.Range("K22").Value = 15.82283202367108
$ws.Range("L22").Value = 10.61026064125042
$ws.Range("N22").Value = 22.64725267302021
$ws.Range("B23").Value = 18.12972212433323
$ws.Range("C23").Value = 5.550899095208331
$ws.Range("E23").Value = 9.512737167400557
$ws.Range("F23").Value = 21.82633154458858
$ws.Range("G23").Value = 3.721013475901505
$ws.Range("I23").Value = 32.5363238143217
$ws.Range("K23").Value = 15.76267729579057
$ws.Range("L23").Value = 10.60008831233058
$ws.Range("N23").Value = 22.66979434703008
$ws.Range("B24").Value = 17.79926858879234
$ws.Range("C24").Value = 5.295832284088608
$ws.Range("E24").Value = 9.503688673066307
$ws.Range("F24").Value = 20.22900810905287
$ws.Range("G24").Value = 3.725253273311026
$ws.Range("I24").Value = 32.56712519263147
$ws.Range("K24").Value = 15.53910199049862
$ws.Range("L24").Value = 10.56485930372254
$ws.Range("N24").Value = 22.75860771712702
$ws.Range("B25").Value = 17.45215718591767
$ws.Range("C25").Value = 5.01203897168057
$ws.Range("E25").Value = 9.499382796996708
$ws.Range("F25").Value = 18.34778573295695
$ws.Range("G25").Value = 3.73015940327246
$ws.Range("I25").Value = 32.62197828381396
$ws.Range("K25").Value = 15.3081025471622
$ws.Range("L25").Value = 10.53400554682491
$ws.Range("N25").Value = 22.86177187419883
